$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110:232 down to 111:233
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new record
$ws.Cells.Item(110, 1).Value = 3
$ws.Cells.Item(110, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(110, 3).Value = "Coquimbo"
$ws.Cells.Item(110, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 9).Date
$ws.Cells.Item(110, 5).Value = 5
$ws.Cells.Item(110, 6).Value = 100112001
$ws.Cells.Item(110, 7).Value = "Berenjena"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 73
$ws.Cells.Item(110, 11).Value = 10000
$ws.Cells.Item(110, 12).Value = 11000
$ws.Cells.Item(110, 13).Value = 10521
$ws.Cells.Item(110, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(110, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(110, 16).Value = 175
$ws.Cells.Item(110, 17).Value = 60
$ws.Cells.Item(110, 18).Value = "Hortaliza"
